# Auto-generated edit script for german_vocab.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vocab")

# --- Fix 5 existing cells: correct missing umlauts / eszett typos ---
$ws.Range("A80").Value = "grün, umweltfreundlich"
$ws.Range("A82").Value = "hilfsbedürftig"
$ws.Range("A99").Value = "um Geschäfte zu machen / fur ihre Geschäfte"
$ws.Range("A125").Value = "Regeln muß man umsetzen / realisieren"
$ws.Range("A140").Value = "Ich wurde eingestellt, um zu helfen, daβ dieser Plan umgesetzt wird."

# --- Append 22 new vocabulary rows (146-167), copying format from row 145 ---
$ws.Range("A145:E145").Copy()
$ws.Range("A146:E167").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A146").Value = "erfolgreich Bankgespräche führen"
$ws.Range("B146").Value = "conducting successful bank conversations"
$ws.Range("C146").Value = 44629
$ws.Range("D146").Value = 6
$ws.Range("E146").Value = "phrase"

$ws.Range("A147").Value = "in welcher Branche ist die Geschäftsidee?"
$ws.Range("B147").Value = "in which sector is the business idea?"
$ws.Range("C147").Value = 44629
$ws.Range("D147").Value = 6
$ws.Range("E147").Value = "phrase"

$ws.Range("A148").Value = "Was hat die Marktanalyse ergeben?"
$ws.Range("B148").Value = "what are the results from the market analysis?"
$ws.Range("C148").Value = 44629
$ws.Range("D148").Value = 6
$ws.Range("E148").Value = "phrase"

$ws.Range("A149").Value = "die Wertvortellung"
$ws.Range("B149").Value = "value proposition"
$ws.Range("C149").Value = 44629
$ws.Range("D149").Value = 6
$ws.Range("E149").Value = "word"

$ws.Range("A150").Value = "das Alleinstellungsmerkmal, die Alleinstellungsmerkmale"
$ws.Range("B150").Value = "the unique selling point"
$ws.Range("C150").Value = 44629
$ws.Range("D150").Value = 6
$ws.Range("E150").Value = "word"

$ws.Range("A151").Value = "der Antragsteller, die Antragsteller"
$ws.Range("B151").Value = "applicant"
$ws.Range("C151").Value = 44629
$ws.Range("D151").Value = 6
$ws.Range("E151").Value = "word"

$ws.Range("A152").Value = "sich darauf vorbereiten, alle Fragen auszudiskutieren"
$ws.Range("B152").Value = "to prepare onself to address all questions"
$ws.Range("C152").Value = 44629
$ws.Range("D152").Value = 6
$ws.Range("E152").Value = "phrase"

$ws.Range("A153").Value = "alle Möglichkeiten parat haben"
$ws.Range("B153").Value = "to have every possibility prepared / covered"
$ws.Range("C153").Value = 44629
$ws.Range("D153").Value = 6
$ws.Range("E153").Value = "phrase"

$ws.Range("A154").Value = "das Risikomanagement sagt, daβ es zu viele Risiken gibt"
$ws.Range("B154").Value = "the risk department says, that there are too many risks"
$ws.Range("C154").Value = 44629
$ws.Range("D154").Value = 6
$ws.Range("E154").Value = "phrase"

$ws.Range("A155").Value = "sie bewerten die Risikograde"
$ws.Range("B155").Value = "they assess the risk level"
$ws.Range("C155").Value = 44629
$ws.Range("D155").Value = 6
$ws.Range("E155").Value = "phrase"

$ws.Range("A156").Value = "der Unternehmer, die Unternehmer"
$ws.Range("B156").Value = "entrepreneur"
$ws.Range("C156").Value = 44629
$ws.Range("D156").Value = 6
$ws.Range("E156").Value = "word"

$ws.Range("A157").Value = "der Geschäftsführer, die Geschäftsführer"
$ws.Range("B157").Value = "general manager"
$ws.Range("C157").Value = 44629
$ws.Range("D157").Value = 6
$ws.Range("E157").Value = "word"

$ws.Range("A158").Value = "den Businessplan erstellen"
$ws.Range("B158").Value = "prepare the business plan"
$ws.Range("C158").Value = 44629
$ws.Range("D158").Value = 6
$ws.Range("E158").Value = "phrase"

$ws.Range("A159").Value = "das Ergebnis der Marktanalyse"
$ws.Range("B159").Value = "the result from the market analysis"
$ws.Range("C159").Value = 44629
$ws.Range("D159").Value = 6
$ws.Range("E159").Value = "phrase"

$ws.Range("A160").Value = "unternehmerische Fähigkeiten"
$ws.Range("B160").Value = "business potential"
$ws.Range("C160").Value = 44629
$ws.Range("D160").Value = 6
$ws.Range("E160").Value = "phrase"

$ws.Range("A161").Value = "die Kredithöhe / die Höhe des Kredits"
$ws.Range("B161").Value = "credit quality"
$ws.Range("C161").Value = 44629
$ws.Range("D161").Value = 6
$ws.Range("E161").Value = "word"

$ws.Range("A162").Value = "die Kreditwürdigkeit"
$ws.Range("B162").Value = "creditworthiness"
$ws.Range("C162").Value = 44629
$ws.Range("D162").Value = 6
$ws.Range("E162").Value = "word"

$ws.Range("A163").Value = "der Blickkontakt - den Blickkontakt halten"
$ws.Range("B163").Value = "eye contact / maintain eye contact"
$ws.Range("C163").Value = 44629
$ws.Range("D163").Value = 6
$ws.Range("E163").Value = "word"

$ws.Range("A164").Value = "das Informationsbedürfnis"
$ws.Range("B164").Value = "information request"
$ws.Range("C164").Value = 44629
$ws.Range("D164").Value = 6
$ws.Range("E164").Value = "word"

$ws.Range("A165").Value = "der Ton"
$ws.Range("B165").Value = "tone"
$ws.Range("C165").Value = 44629
$ws.Range("D165").Value = 6
$ws.Range("E165").Value = "word"

$ws.Range("A166").Value = "die Nachfrage (von)"
$ws.Range("B166").Value = "demand (as in supply and demand)"
$ws.Range("C166").Value = 44629
$ws.Range("D166").Value = 6
$ws.Range("E166").Value = "word"

$ws.Range("A167").Value = "ruhig / gelassen bleiben"
$ws.Range("B167").Value = "stay calm"
$ws.Range("C167").Value = 44629
$ws.Range("D167").Value = 6
$ws.Range("E167").Value = "word"

# --- Update view state to match the target selection/scroll position ---
$ws.Activate()
$ws.Range("A161").Select()
$excel.ActiveWindow.ScrollRow = 135
$excel.ActiveWindow.ScrollColumn = 1

